$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the header label in B2: it used to be the stray pandas artifact
# "unnamed: 1_level_1" and should read "total" (matching B1's pair column).
$ws.Range("B2").Value = "total"

# Two spacer/header rows ("situação do domicílio" and "grandes regiões e
# unidades da federação") were removed from the data table; deleting them
# shifts every following data row up so each region/age-group label stays
# aligned with its own numeric values.
$ws.Rows("5").Delete()
$ws.Rows("7").Delete()
